$d = $word.ActiveDocument

# --------------------------------------------------------------------
# Helper: place a zero-width bookmark at the very end of a paragraph's
# text (i.e. right after the last run, right before the paragraph
# mark). Doing this directly by collapsing a Range to that exact
# boundary position is unreliable here, so work around it: insert a
# one-character pad right before the paragraph mark, bookmark that
# (safely non-boundary) character, then delete the pad again -- the
# bookmark collapses back down to the correct end-of-text position.
# Giving the bookmark the reserved name "_GoBack" also removes any
# previous bookmark of that name elsewhere in the document, since
# bookmark names must be unique -- which is exactly the "move the
# _GoBack bookmark" behaviour we want.
# --------------------------------------------------------------------
function Set-BookmarkAtParagraphEnd($paraIndex, $bookmarkName) {
    $p = $d.Paragraphs.Item($paraIndex)
    $r = $p.Range.Duplicate
    $r.MoveEnd(1, -1)
    $r.Collapse(0)
    $r.InsertAfter("@")

    $p2 = $d.Paragraphs.Item($paraIndex)
    $full = $p2.Range.Duplicate
    $full.MoveEnd(1, -1)
    $padChar = $full.Duplicate
    $padChar.MoveStart(1, ($full.End - $full.Start - 1))
    $d.Bookmarks.Add($bookmarkName, $padChar)

    $p3 = $d.Paragraphs.Item($paraIndex)
    $full2 = $p3.Range.Duplicate
    $full2.MoveEnd(1, -1)
    $padChar2 = $full2.Duplicate
    $padChar2.MoveStart(1, ($full2.End - $full2.Start - 1))
    $padChar2.Text = ""
}

# The document contains a duplicated "Geen validatie" stub section
# (heading + empty Beschrijving/Oplossing bodies) right after the "Knop
# ... botst met andere items" section. That stub gets turned into real
# content describing a "Strings niet in de strings.xml" bug/fix, and
# the _GoBack bookmark (previously at the end of the preceding
# paragraph) moves to the end of the new final paragraph of this
# section.

# Paragraph 21: Kop2 "Geen validatie" stub heading -> new heading text
$d.Paragraphs.Item(21).Range.Text = "Strings niet in de strings.xml"

# Paragraph 22: Kop3 "Beschrijving" -- unchanged, already correct.

# Paragraph 23: empty body under "Beschrijving" -> filled in
$d.Paragraphs.Item(23).Range.Text = "De strings werden hardcoded toegevoegd in plaats via de strings.xml"

# Paragraph 24: Kop3 "Oplossing" -- unchanged, already correct.

# Paragraph 25: empty body under "Oplossing" -> filled in
$d.Paragraphs.Item(25).Range.Text = "De strings in de strings.xml gezet"

# Move the _GoBack bookmark to the end of the text just inserted above.
Set-BookmarkAtParagraphEnd 25 "_GoBack"
